$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.382.15'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.914.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.97%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.730'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '256.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.28%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.70'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.369'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.60%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0765'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0990'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.190.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.738'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.921.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.368.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '75.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0853'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '246.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.38%  '

$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("E25").Value = '  +7.15%  '

$ws.Range("E26").Value = '  +2.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.15%  '

$ws.Range("E30").Value = '  +5.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.127.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +24.80%  '

$ws.Range("E34").Value = '  +14.26%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0594'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.38%  '

$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.923'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.82%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.79%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0647'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.49'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.348.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.12%  '

$ws.Range("E47").Value = '  +1.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.59%  '

$ws.Range("E49").Value = '  -0.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '45.50'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.39%  '

$ws.Range("E51").Value = '  +6.88%  '
